$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new rows of test data matching existing pattern (Function col B, Test col C)
$ws.Range("B9").Value = "displaySaveCalibYN"
$ws.Range("C9").Value = """"

$ws.Range("B10").Value = "displayNewCalibYN"
$ws.Range("C10").Value = """"

# Match the styling of the existing "Status" (D) column cells (D2:D8 style)
$ws.Range("D8").Copy()
$ws.Range("D9:D10").PasteSpecial(-4122)

# Update the active selection to match the new working cell
$ws.Range("F8").Select()
